$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are stored as text, matching the
# original inline-string cell contents (many prices are plain decimals that
# Excel would otherwise auto-convert to numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.287.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.584.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.807.41"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.585.71"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.298.40"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.26"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.02%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.30"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.25%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.65%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.283.92"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.613"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.54%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.45"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.719.83"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.85"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.00%  "

# Row 50: coin replaced (BabyDogeCoin -> EnergySwap) with new link, price, volume
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.47"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.28%  "
